# Siemens sheet: add a "weight" column between "order" and "frequency",
# giving each row a sequential weight (1..10), and keep the "frequency"
# column (all 100s) shifted one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F; the existing "frequency" column (F)
# shifts right to become G.
$ws.Columns.Item(6).Insert()

# New header for column F.
$ws.Cells.Item(1, 6).Value = "weight"

# Sequential weight values for rows 2..11 (1,2,3,...,10).
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 6).Value = $r - 1
}

# Match column F's width to the rest of that block (same as column G).
$ws.Range("F1:G11").ColumnWidth = 13.14

# Column A got narrower in this revision.
$ws.Columns.Item(1).ColumnWidth = 17.43

# Reflect the edited area in the current selection, as the author left it.
$ws.Range("F2:F11").Select() | Out-Null
